{"js": "// Map of old answer-line text -> new answer-line text, taken from the\n// authoritative OOXML diff. Each table cell holds a single run whose\n// text is a unique \"A\u00d7B=C\" string, so searching the whole body for that\n// exact string and replacing it is unambiguous for every entry.\nconst replacements = [\n  [\"67\u00d755=3685\", \"67\u00d779=5293\"],\n  [\"68\u00d740=2720\", \"98\u00d724=2352\"],\n  [\"56\u00d762=3472\", \"30\u00d751=1530\"],\n  [\"47\u00d725=1175\", \"66\u00d730=1980\"],\n  [\"57\u00d789=5073\", \"84\u00d791=7644\"],\n  [\"46\u00d730=1380\", \"53\u00d781=4293\"],\n  [\"30\u00d764=1920\", \"32\u00d733=1056\"],\n  [\"34\u00d738=1292\", \"60\u00d743=2580\"],\n  [\"38\u00d782=3116\", \"63\u00d769=4347\"],\n  [\"57\u00d795=5415\", \"74\u00d738=2812\"],\n  [\"16\u00d787=1392\", \"79\u00d757=4503\"],\n  [\"98\u00d756=5488\", \"26\u00d713=338\"],\n  [\"93\u00d752=4836\", \"46\u00d738=1748\"],\n  [\"88\u00d757=5016\", \"93\u00d718=1674\"],\n  [\"51\u00d712=612\", \"19\u00d758=1102\"],\n  [\"91\u00d789=8099\", \"59\u00d795=5605\"],\n  [\"49\u00d723=1127\", \"31\u00d792=2852\"],\n  [\"42\u00d733=1386\", \"64\u00d725=1600\"],\n  [\"64\u00d733=2112\", \"53\u00d729=1537\"],\n  [\"34\u00d759=2006\", \"54\u00d763=3402\"],\n  [\"74\u00d787=6438\", \"95\u00d747=4465\"],\n  [\"53\u00d756=2968\", \"94\u00d724=2256\"],\n  [\"98\u00d729=2842\", \"17\u00d796=1632\"],\n  [\"11\u00d711=121\", \"88\u00d785=7480\"],\n  [\"19\u00d734=646\", \"58\u00d739=2262\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${oldText}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Map of old answer-line text -> new answer-line text, taken from the\n# authoritative OOXML diff. Each <w:t> run in the table is a unique\n# \"A\u00d7B=C\" string, so a plain text Find/Replace on the whole document\n# body is unambiguous for every entry.\n$replacements = @(\n    @{ Old = '67\u00d755=3685'; New = '67\u00d779=5293' }\n    @{ Old = '68\u00d740=2720'; New = '98\u00d724=2352' }\n    @{ Old = '56\u00d762=3472'; New = '30\u00d751=1530' }\n    @{ Old = '47\u00d725=1175'; New = '66\u00d730=1980' }\n    @{ Old = '57\u00d789=5073'; New = '84\u00d791=7644' }\n    @{ Old = '46\u00d730=1380'; New = '53\u00d781=4293' }\n    @{ Old = '30\u00d764=1920'; New = '32\u00d733=1056' }\n    @{ Old = '34\u00d738=1292'; New = '60\u00d743=2580' }\n    @{ Old = '38\u00d782=3116'; New = '63\u00d769=4347' }\n    @{ Old = '57\u00d795=5415'; New = '74\u00d738=2812' }\n    @{ Old = '16\u00d787=1392'; New = '79\u00d757=4503' }\n    @{ Old = '98\u00d756=5488'; New = '26\u00d713=338' }\n    @{ Old = '93\u00d752=4836'; New = '46\u00d738=1748' }\n    @{ Old = '88\u00d757=5016'; New = '93\u00d718=1674' }\n    @{ Old = '51\u00d712=612'; New = '19\u00d758=1102' }\n    @{ Old = '91\u00d789=8099'; New = '59\u00d795=5605' }\n    @{ Old = '49\u00d723=1127'; New = '31\u00d792=2852' }\n    @{ Old = '42\u00d733=1386'; New = '64\u00d725=1600' }\n    @{ Old = '64\u00d733=2112'; New = '53\u00d729=1537' }\n    @{ Old = '34\u00d759=2006'; New = '54\u00d763=3402' }\n    @{ Old = '74\u00d787=6438'; New = '95\u00d747=4465' }\n    @{ Old = '53\u00d756=2968'; New = '94\u00d724=2256' }\n    @{ Old = '98\u00d729=2842'; New = '17\u00d796=1632' }\n    @{ Old = '11\u00d711=121'; New = '88\u00d785=7480' }\n    @{ Old = '19\u00d734=646'; New = '58\u00d739=2262' }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $true,       # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\nWrite-Output \"Replaced $($replacements.Count) answer cells\"\n"}
